# ---------------------------------------------------------------------------
# Adds the 2022-Q4 quarter to the 688029-南微医学 workbook:
#   1. A new worksheet "2022-Q4" is inserted right after "总计" (and before
#      "2022-Q3"), carrying the per-fund holding detail for that quarter.
#   2. The "总计" (summary) sheet gets a new top data row for 2022-Q4 and all
#      the previously existing quarters shift down by one row.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" summary sheet: insert the 2022-Q4 row, push the rest
#    down, and land 2020-Q4 (previously the last row) in the new last row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Shift existing B/C/D data down by one row (bottom-up so we never clobber a
# value before it has been copied). Column A already holds the sequential
# 0,1,2,... index for every row and does not need to move.
for ($r = 9; $r -ge 2; $r--) {
    $dst = $r + 1
    $total.Cells.Item($dst, 2).Value = $total.Cells.Item($r, 2).Value2
    $total.Cells.Item($dst, 3).Value = $total.Cells.Item($r, 3).Value2
    $total.Cells.Item($dst, 4).Value = $total.Cells.Item($r, 4).Value2
}

# New last row (row 10) needs the sequential index too.
$total.Cells.Item(10, 1).Value = 8
$srcIndexCell = $total.Cells.Item(9, 1)
$srcIndexCell.Copy()
$total.Cells.Item(10, 1).PasteSpecial(-4122)

# Fresh 2022-Q4 figures go into row 2 (index 0 is already there).
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 13
$total.Cells.Item(2, 4).Value = 1.31

# ---------------------------------------------------------------------------
# 2. Insert the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------------------
$headerSrc = $total.Cells.Item(1, 2)

$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q4.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

$data = @(
    @(0,  "000893", "工银创新动力股票",                 "15.58", "81.84", "3.40", "0.5297", 8),
    @(1,  "506007", "广发科创板两年定开混合",             "5.31",  "88.81", "5.08", "0.2697", 4),
    @(2,  "010405", "惠升医药健康6个月持有期混合",         "9.98",  "72.67", "2.54", "0.2535", 8),
    @(3,  "001695", "泓德泓业灵活配置混合",               "1.56",  "90.00", "4.47", "0.0697", 6),
    @(4,  "002801", "泓德泓信灵活配置混合",               "1.75",  "92.52", "3.35", "0.0586", 9),
    @(5,  "011781", "泓德慧享混合A",                     "5.92",  "27.55", "0.96", "0.0568", 6),
    @(6,  "009015", "泓德睿享一年持有期混合A",             "2.52",  "28.92", "1.29", "0.0325", 6),
    @(7,  "002681", "金鹰元和灵活配置混合A",               "0.30",  "81.19", "6.09", "0.0183", 1),
    @(8,  "002682", "金鹰元和灵活配置混合C",               "0.23",  "81.19", "6.09", "0.0140", 1),
    @(9,  "350008", "天治新消费灵活配置混合",             "0.12",  "94.01", "3.24", "0.0039", 5),
    @(10, "515590", "前海开源中证500等权重ETF",           "0.34",  "95.19", "0.29", "0.0010", 7),
    @(11, "009016", "泓德睿享一年持有期混合C",             "0.06",  "28.92", "1.29", "0.0008", 6),
    @(12, "011782", "泓德慧享混合C",                     "0.00",  "27.55", "0.96", $null,    6)
)

foreach ($row in $data) {
    $r = 2 + $row[0]
    $q4.Cells.Item($r, 1).Value = $row[0]

    $q4.Cells.Item($r, 2).NumberFormat = "@"
    $q4.Cells.Item($r, 2).Value = $row[1]

    $q4.Cells.Item($r, 3).Value = $row[2]

    $q4.Cells.Item($r, 4).NumberFormat = "@"
    $q4.Cells.Item($r, 4).Value = $row[3]

    $q4.Cells.Item($r, 5).NumberFormat = "@"
    $q4.Cells.Item($r, 5).Value = $row[4]

    $q4.Cells.Item($r, 6).NumberFormat = "@"
    $q4.Cells.Item($r, 6).Value = $row[5]

    if ($row[6] -eq $null) {
        $q4.Cells.Item($r, 7).Value = 0
    } else {
        $q4.Cells.Item($r, 7).NumberFormat = "@"
        $q4.Cells.Item($r, 7).Value = $row[6]
    }

    $q4.Cells.Item($r, 8).Value = $row[7]
}

# Header row + the index column share the bold/centered/bordered look used
# throughout the workbook ("style 2"). Copy that formatting across in one go.
$headerSrc.Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

$indexSrc = $total.Cells.Item(2, 1)
$indexSrc.Copy()
$q4.Range("A2:A14").PasteSpecial(-4122)
